$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("AR1").Value = 0.92957696593215666
$ws.Range("BP1").Value = 0.78415161858768401
$ws.Range("BM2").Value = 0.79853757683489657
$ws.Range("AN3").Value = 0.86675297542639285
$ws.Range("AT3").Value = 0.98803741519232324
$ws.Range("BM3").Value = 0.97467084394176418
$ws.Range("F4").Value = 0.62625711015017771
$ws.Range("R4").Value = 0.93814534976730712
$ws.Range("O5").Value = 0.79790238709322847
$ws.Range("F8").Value = 0.67125795954864675
$ws.Range("T9").Value = 0.94887493783895138
$ws.Range("V9").Value = 0.9184915094031989
$ws.Range("AP9").Value = 0.96937334373060569
$ws.Range("L10").Value = 0.7804974477554345
$ws.Range("AW10").Value = 0.90351687061132857
$ws.Range("BL10").Value = 0.99380831427699423
$ws.Range("W11").Value = 0.8584548301996503
$ws.Range("BA11").Value = 0.94364272360725843
$ws.Range("BH11").Value = 0.61846329666006095
$ws.Range("H12").Value = 0.78633029028066437
$ws.Range("O13").Value = 0.81810341218615235
$ws.Range("AP13").Value = 0.85271560276837433
$ws.Range("AX13").Value = 0.99443709697710614
$ws.Range("F15").Value = 0.63554748137020889
$ws.Range("AO15").Value = 0.93093343431451214
$ws.Range("G16").Value = 0.79216742148820452
$ws.Range("N16").Value = 0.70554012816253175
$ws.Range("AD16").Value = 0.92125401607148005
$ws.Range("AV16").Value = 0.93464108112162347
$ws.Range("BP16").Value = 0.70332012752060646
$ws.Range("AW17").Value = 0.78665618886472899
$ws.Range("AX17").Value = 0.85243278671776423
$ws.Range("BC17").Value = 0.72696103351835983
$ws.Range("O18").Value = 0.96442288401593979
$ws.Range("E19").Value = 0.8812038023523332
$ws.Range("E20").Value = 0.8081945529964728
$ws.Range("K20").Value = 0.73786341252292575
$ws.Range("B21").Value = 0.80161034449233481
$ws.Range("F21").Value = 0.94550258220418493
$ws.Range("U23").Value = 0.78143366014060045
$ws.Range("W24").Value = 0.96792828592148794
$ws.Range("AA24").Value = 0.99064502157643752
$ws.Range("AU24").Value = 0.95204052607950984
$ws.Range("X25").Value = 0.75823588238105089
$ws.Range("Z25").Value = 0.58352569973355894
$ws.Range("R26").Value = 0.71921743292993034
$ws.Range("AK26").Value = 0.53219694260678541
$ws.Range("AL26").Value = 0.86093327095808858
$ws.Range("BI26").Value = 0.763749341733025
$ws.Range("D27").Value = 0.63871500336188536
$ws.Range("AH27").Value = 0.90030566683102342
$ws.Range("AN27").Value = 0.56400163972864048
$ws.Range("BG27").Value = 0.80549310249529693
$ws.Range("Z28").Value = 0.98907036341755183
$ws.Range("BK28").Value = 0.74667586546891307
$ws.Range("V29").Value = 0.93246506051904876
$ws.Range("AS29").Value = 0.57713401774223838
$ws.Range("BI29").Value = 0.9512919262844155
$ws.Range("C30").Value = 0.5641226821529759
$ws.Range("AX30").Value = 0.93714150089371318
$ws.Range("BJ30").Value = 0.67918673149659148
$ws.Range("X31").Value = 0.99316351686560533
$ws.Range("AG31").Value = 0.62478521645952989
$ws.Range("BD31").Value = 0.89542314960633562
$ws.Range("AH32").Value = 0.86925871670030674
$ws.Range("Q33").Value = 0.92464007425748651
$ws.Range("AF33").Value = 0.9714749127929011
$ws.Range("BM33").Value = 0.76898309791059605
$ws.Range("AJ34").Value = 0.85708912346301913
$ws.Range("AO34").Value = 0.96745004343059193
$ws.Range("AV34").Value = 0.79893800583716346
$ws.Range("L35").Value = 0.92034973780310958
$ws.Range("V35").Value = 0.78698267753933293
$ws.Range("G36").Value = 0.86451821944838692
$ws.Range("AC36").Value = 0.92230778853341566
$ws.Range("AL36").Value = 0.95566153528034614
$ws.Range("AI37").Value = 0.94490844807819974
$ws.Range("AM37").Value = 0.96039578962858074
$ws.Range("BL37").Value = 0.95469589858556869
$ws.Range("AY38").Value = 0.87711670518517171
$ws.Range("A39").Value = 0.73438722394983769
$ws.Range("AO39").Value = 0.77924633052378356
$ws.Range("BA39").Value = 0.78018175717117044
$ws.Range("BJ39").Value = 0.89895365614754441
$ws.Range("AQ40").Value = 0.60035967194625783
$ws.Range("AF41").Value = 0.85421023595547385
$ws.Range("AF42").Value = 0.85383413412198728
$ws.Range("C43").Value = 0.60055808306234315
$ws.Range("AG43").Value = 0.62217877659435206
$ws.Range("AO43").Value = 0.81616974204339821
$ws.Range("J44").Value = 0.98706228167311627
$ws.Range("O44").Value = 0.94425092197857163
$ws.Range("BA44").Value = 0.79991365276062676
$ws.Range("N45").Value = 0.68425908034559924
$ws.Range("H46").Value = 0.93944912814751014
$ws.Range("AR46").Value = 0.88925880900406895
$ws.Range("BL46").Value = 0.98910938923504932
$ws.Range("T47").Value = 0.6829617407261368
$ws.Range("AR47").Value = 0.99906141342462162
$ws.Range("BO47").Value = 0.88256521042121849
$ws.Range("AY48").Value = 0.89318838048744531
$ws.Range("Y49").Value = 0.83700566637854901
$ws.Range("BE49").Value = 0.85743041297769906
$ws.Range("AF50").Value = 0.76323498931429345
$ws.Range("BG50").Value = 0.69608135863357012
$ws.Range("H53").Value = 0.84607568586839565
$ws.Range("AZ53").Value = 0.90125540160170925
$ws.Range("BG53").Value = 0.57481459062163165
$ws.Range("AZ54").Value = 0.68447428745668748
$ws.Range("BA54").Value = 0.8848658251180066
$ws.Range("K55").Value = 0.8792376391997343
$ws.Range("S55").Value = 0.87762460635670947
$ws.Range("BA55").Value = 0.99816421335729788
$ws.Range("M56").Value = 0.73380930213051587
$ws.Range("N56").Value = 0.99897986133339867
$ws.Range("AU57").Value = 0.7857651722357093
$ws.Range("BC57").Value = 0.59143694828679627
$ws.Range("X58").Value = 0.99382643351201017
$ws.Range("AT58").Value = 0.58957314743178069
$ws.Range("BE58").Value = 0.73818140061930837
$ws.Range("BH58").Value = 0.92607688347465955
$ws.Range("AG60").Value = 0.95584679951220786
$ws.Range("AV60").Value = 0.89054247189628732
$ws.Range("AG61").Value = 0.90015189268249918
$ws.Range("AL61").Value = 0.94099577015708324
$ws.Range("BJ61").Value = 0.95052813902697197
$ws.Range("E62").Value = 0.9608511588163362
$ws.Range("BC63").Value = 0.9290631172168804
$ws.Range("BN63").Value = 0.76621599302940924
$ws.Range("AC64").Value = 0.85256331881808767
$ws.Range("BP64").Value = 0.83212255611730246
$ws.Range("AB65").Value = 0.82973514069287901
$ws.Range("BB65").Value = 0.86613068728913323
$ws.Range("I66").Value = 0.97478409071597116
$ws.Range("BN67").Value = 0.89530472804725292
$ws.Range("BP67").Value = 0.90331362098523171
